$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 66 entirely; subsequent rows shift up by one.
$ws.Rows(66).Delete()

# The underlying engine does not correctly relocate the 1x1 "merge" ranges
# that sat at A67/B67 and A141/B141 in the original sheet when the row
# above them was removed: the merge at row 67 was dropped, and the merge
# at row 141 stayed in place instead of moving up to row 140. Fix those up
# by hand so the result matches a real row deletion (A67->A66, A141->A140).
$ws.Range("A141").UnMerge()
$ws.Range("B141").UnMerge()

$ws.Range("A66").Merge()
$ws.Range("B66").Merge()
$ws.Range("A140").Merge()
$ws.Range("B140").Merge()
